$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the anchor paragraph: the blank paragraph 5 positions after the
# "the-normal-matrix" link paragraph (end of section 9 / start of the large
# block of trailing blank paragraphs at the end of the document).
# ---------------------------------------------------------------------------
$anchorIndex = 126
$anchorPara = $d.Paragraphs.Item($anchorIndex)
$insertRange = $anchorPara.Range
$insertRange.Collapse(1)

# ---------------------------------------------------------------------------
# Insert the new "10. GLSL function Specification" heading paragraph plus
# the accompanying URL paragraph right before that anchor paragraph.
# ---------------------------------------------------------------------------
$newXml = '<?xml version="1.0" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p>' +
  '<w:pPr>' +
  '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' +
  '<w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>' +
  '</w:pPr>' +
  '<w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>10</w:t></w:r>' +
  '<w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r>' +
  '<w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>GLSL function Specification</w:t></w:r>' +
  '</w:p>' +
  '<w:p>' +
  '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr>' +
  '<w:r><w:t>http://www.khronos.org/registry/gles/specs/2.0/GLSL_ES_Specification_1.0.17.pdf</w:t></w:r>' +
  '</w:p>' +
  '</w:body>' +
  '</w:document>' +
  '</pkg:xmlData>' +
  '</pkg:part>' +
  '</pkg:package>'

$insertRange.InsertXML($newXml)

# ---------------------------------------------------------------------------
# 55 blank paragraphs after the new content are left untouched. After them,
# a run of 8 blank paragraphs collapses down to a single blank paragraph
# that gains an eastAsia font hint on its paragraph mark.
# ---------------------------------------------------------------------------
$collapseStart = $anchorIndex + 55
$collapseEnd = $collapseStart + 8 - 1

$startRange = $d.Paragraphs.Item($collapseStart).Range
$endRange = $d.Paragraphs.Item($collapseEnd).Range
$collapseRange = $d.Range($startRange.Start, $endRange.End)

$hintXml = '<?xml version="1.0" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p>' +
  '<w:pPr>' +
  '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' +
  '<w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr>' +
  '</w:pPr>' +
  '</w:p>' +
  '</w:body>' +
  '</w:document>' +
  '</pkg:xmlData>' +
  '</pkg:part>' +
  '</pkg:package>'

$collapseRange.InsertXML($hintXml)

Write-Host "Final paragraph count:" $d.Paragraphs.Count
